$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header
$ws.Range("A1").Value = "URL"

# Row 2: the bypass URL, now reachable without the :443 port, as a hyperlink
$target = "https://pins-test.gopro.net/selfservice/web/portal/exemption.html"
$ws.Hyperlinks.Add($ws.Range("A2"), $target) | Out-Null

$ws.Range("A2").Select() | Out-Null
